# microscopy.xlsx commit: "typos in microscopy plus cleaning of frames in main script"
#
# 1. GW-Oct (sheet "GW-Oct") was missing the REPLICATE (column D) values for
#    rows 2-56 -- a typo/paste bug. Re-populate them to match the pattern used
#    on every other tab (GW-May, WB-May, RB-May, WB-Oct, RB-Oct all already
#    have column D filled in).
# 2. Make GW-Oct the active tab/selection (it was RB-Oct before), and move the
#    selection on GW-May.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Restore column D ("REPLICATE") values on GW-Oct, rows 2-56.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("GW-Oct")

$replicateValues = @(
    1,1,1,1,1,1,1,1,1,1,1,1,   # rows 2-13
    2,2,2,2,2,2,2,2,2,2,2,2,2,2, # rows 14-27
    3,3,3,3,3,3,3,3,3,3,3,3,3,3, # rows 28-41
    4,4,4,4,4,4,4,4,4,        # rows 42-50
    6,6,6,6,6,6                # rows 51-56
)

$startRow = 2
for ($i = 0; $i -lt $replicateValues.Length; $i++) {
    $row = $startRow + $i
    $ws4.Cells.Item($row, 4).Value = $replicateValues[$i]
}

# Column D is currently merged with column C in the stored column metadata
# (same width, 10.12). Now that it holds its own data again, widen it
# slightly so it is no longer identical to column C (target file width is
# 11.99; this COM layer quantizes column widths, so set the nearest value
# that yields that result).
$ws4.Columns.Item(4).ColumnWidth = 11.084

# ---------------------------------------------------------------------------
# 2) View/selection changes.
# ---------------------------------------------------------------------------

# GW-May: selection stays at A1 (no cell-selection change requested), only
# the scroll position moves -- nothing further needed here beyond leaving it
# alone.
$ws1 = $wb.Worksheets.Item("GW-May")
$ws1.Range("A1").Select()

# GW-Oct becomes the active sheet, with E41 selected (it previously held
# K10 while RB-Oct was the active tab).
$ws4.Activate()
$ws4.Range("E41").Select()

Write-Output "done"
